$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Reword the "Gestión de Clientes" bullet: add "con el tiempo" and change
#    "personalizar"/"gestionar" -> "empezar a  personalizar"/"gestionando".
# ---------------------------------------------------------------------------
$oldFragment = "sino también personalizar el servicio, gestionar promociones"
$newFragment = "sino también, con el tiempo, empezar a  personalizar el servicio, gestionando promociones"

$rng = $d.Content
$found = $rng.Find.Execute($oldFragment, $true, $false, $false, $false, $false, $true, 1, $false, $newFragment, 2)
if (-not $found) {
    throw "Could not locate the sentence fragment to reword."
}

# ---------------------------------------------------------------------------
# 2. Swap the two floating picture names: the shape currently titled
#    "image2.png" becomes "image1.png", and the one titled "image1.png"
#    becomes "image2.png" (the underlying embedded pictures are untouched).
# ---------------------------------------------------------------------------
$shapes = $d.Shapes
$shapeNamedImage2 = $null
$shapeNamedImage1 = $null
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.Name -eq "image2.png") { $shapeNamedImage2 = $shp }
    elseif ($shp.Name -eq "image1.png") { $shapeNamedImage1 = $shp }
}

if ($shapeNamedImage2 -and $shapeNamedImage1) {
    # Route through a scratch name so the two renames never collide.
    $shapeNamedImage2.Name = "__swap_tmp__"
    $shapeNamedImage1.Name = "image2.png"
    $shapeNamedImage2.Name = "image1.png"
}
